$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 29.80827733333334
$ws.Range("N2").Value = 89.42483200000001
$ws.Range("O2").Value = 0.4866975737940222
$ws.Range("P2").Value = 0.4866975737940221
$ws.Range("Q2").Value = 36.967599574976
$ws.Range("R2").Value = 332.708396174784
$ws.Range("S2").Value = 0.4866975737940222
$ws.Range("T2").Value = 0.4866975737940221

# Row 3
$ws.Range("O3").Value = 0.437868100938039
$ws.Range("P3").Value = 0.437868100938039
$ws.Range("S3").Value = 0.437868100938039
$ws.Range("T3").Value = 0.437868100938039

# Row 4
$ws.Range("M4").Value = 4.620050333333332
$ws.Range("O4").Value = 0.07543432526793886
$ws.Range("P4").Value = 0.07543432526793886
$ws.Range("S4").Value = 0.07543432526793886
$ws.Range("T4").Value = 0.07543432526793886
